$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF"), matching the style used by the
# rest of the header row (bold + border + centered/top alignment = same
# style as H1).
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Numeric data for the new I0 / IF columns, rows 2-24.
$data = @{
    2  = @(1, 2)
    3  = @(9, 9)
    4  = @(9, 9)
    5  = @(5, 5)
    6  = @(9, 9)
    7  = @(2, 2)
    8  = @(2, 4)
    9  = @(5, 5)
    10 = @(4, 5)
    11 = @(8, 9)
    12 = @(9, 9)
    13 = @(7, 7)
    14 = @(5, 6)
    15 = @(8, 8)
    16 = @(6, 7)
    17 = @(7, 8)
    18 = @(7, 7)
    19 = @(3, 7)
    20 = @(4, 4)
    21 = @(5, 7)
    22 = @(6, 8)
    23 = @(7, 8)
    24 = @(4, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
